$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$errorText = @'
Checked availability: Failed to select the date: Message: no such element: Unable to locate element: {"method":"css selector","selector":"#restProfileSideBarDtpDayPicker-label"}
  (Session info: chrome=129.0.6668.71); For documentation on this error, please visit: https://www.selenium.dev/documentation/webdriver/troubleshooting/errors#no-such-element-exception
Stacktrace:
	GetHandleVerifier [0x00007FF60823B645+29573]
	(No symbol) [0x00007FF6081B0470]
	(No symbol) [0x00007FF60806B6EA]
	(No symbol) [0x00007FF6080BF815]
	(No symbol) [0x00007FF6080BFA6C]
	(No symbol) [0x00007FF60810B917]
	(No symbol) [0x00007FF6080E733F]
	(No symbol) [0x00007FF6081086BC]
	(No symbol) [0x00007FF6080E70A3]
	(No symbol) [0x00007FF6080B12DF]
	(No symbol) [0x00007FF6080B2441]
	GetHandleVerifier [0x00007FF60856C58D+3375821]
	GetHandleVerifier [0x00007FF6085B7987+3684039]
	GetHandleVerifier [0x00007FF6085ACDAB+3640043]
	GetHandleVerifier [0x00007FF6082FB7C6+816390]
	(No symbol) [0x00007FF6081BB77F]
	(No symbol) [0x00007FF6081B75A4]
	(No symbol) [0x00007FF6081B7740]
	(No symbol) [0x00007FF6081A659F]
	BaseThreadInitThunk [0x00007FFF5DA8257D+29]
	RtlUserThreadStart [0x00007FFF5ECEAF08+40]

'@

$ws.Range('A240').Value = @'
2024-10-01 20:33:16
'@
$ws.Range('B240').Value = @'
check_availability
'@
$ws.Range('C240').Value = @'
https://www.bestbuy.com/site/microsoft-xbox-wireless-controller-for-xbox-series-x-xbox-series-s-xbox-one-windows-devices-sky-cipher-special-edition/6584960.p?skuId=6584960
'@
$ws.Range('D240').NumberFormat = '@'
$ws.Range('D240').Value = @'
$69.99
'@
$ws.Range('D240').Style = 'Normal'
$ws.Range('E240').NumberFormat = '@'
$ws.Range('E240').Value = @'
2024-10-01
'@
$ws.Range('E240').Style = 'Normal'
$ws.Range('F240').Value = @'
20:33:16
'@

$ws.Range('A241').Value = @'
2024-10-01 20:33:22
'@
$ws.Range('B241').Value = @'
check_availability
'@
$ws.Range('C241').Value = @'
https://www.bestbuy.com/site/microsoft-xbox-wireless-controller-for-xbox-series-x-xbox-series-s-xbox-one-windows-devices-sky-cipher-special-edition/6584960.p?skuId=6584960
'@
$ws.Range('D241').NumberFormat = '@'
$ws.Range('D241').Value = @'
$69.99
'@
$ws.Range('D241').Style = 'Normal'
$ws.Range('E241').NumberFormat = '@'
$ws.Range('E241').Value = @'
2024-10-01
'@
$ws.Range('E241').Style = 'Normal'
$ws.Range('F241').Value = @'
20:33:22
'@

$ws.Range('A242').Value = @'
2024-10-01 20:33:59
'@
$ws.Range('B242').Value = @'
check_availability
'@
$ws.Range('C242').Value = @'
https://www.opentable.com/r/the-rux-nashville
'@
$ws.Range('D242').Value = @'
Checked availability: Selected or default date current date is available for booking.
'@
$ws.Range('E242').NumberFormat = '@'
$ws.Range('E242').Value = @'
2024-10-01
'@
$ws.Range('E242').Style = 'Normal'
$ws.Range('F242').Value = @'
20:33:59
'@

$ws.Range('A243').Value = @'
2024-10-01 20:34:04
'@
$ws.Range('B243').Value = @'
check_availability
'@
$ws.Range('C243').Value = @'
https://www.opentable.com/r/hals-the-steakhouse-nashville
'@
$ws.Range('D243').Value = $errorText
$ws.Range('E243').NumberFormat = '@'
$ws.Range('E243').Value = @'
2024-10-01
'@
$ws.Range('E243').Style = 'Normal'
$ws.Range('F243').Value = @'
20:34:04
'@

$ws.Range('A244').Value = @'
2024-10-01 20:34:23
'@
$ws.Range('B244').Value = @'
check_availability
'@
$ws.Range('C244').Value = @'
https://www.opentable.com/r/hals-the-steakhouse-nashville
'@
$ws.Range('D244').Value = $errorText
$ws.Range('E244').NumberFormat = '@'
$ws.Range('E244').Value = @'
2024-10-01
'@
$ws.Range('E244').Style = 'Normal'
$ws.Range('F244').Value = @'
20:34:23
'@

$ws.Range('A245').Value = @'
2024-10-01 20:34:41
'@
$ws.Range('B245').Value = @'
check_availability
'@
$ws.Range('C245').Value = @'
https://www.opentable.com/r/hals-the-steakhouse-nashville
'@
$ws.Range('D245').Value = $errorText
$ws.Range('E245').NumberFormat = '@'
$ws.Range('E245').Value = @'
2024-10-01
'@
$ws.Range('E245').Style = 'Normal'
$ws.Range('F245').Value = @'
20:34:41
'@

$ws.Range('A246').Value = @'
2024-10-01 20:35:06
'@
$ws.Range('B246').Value = @'
check_availability
'@
$ws.Range('C246').Value = @'
https://www.opentable.com/r/sinatra-bar-and-lounge-nashville
'@
$ws.Range('D246').Value = $errorText
$ws.Range('E246').NumberFormat = '@'
$ws.Range('E246').Value = @'
2024-10-01
'@
$ws.Range('E246').Style = 'Normal'
$ws.Range('F246').Value = @'
20:35:06
'@

$ws.Range('A247').Value = @'
2024-10-01 20:35:25
'@
$ws.Range('B247').Value = @'
check_availability
'@
$ws.Range('C247').Value = @'
https://www.opentable.com/r/sinatra-bar-and-lounge-nashville
'@
$ws.Range('D247').Value = $errorText
$ws.Range('E247').NumberFormat = '@'
$ws.Range('E247').Value = @'
2024-10-01
'@
$ws.Range('E247').Style = 'Normal'
$ws.Range('F247').Value = @'
20:35:25
'@

$ws.Range('A248').Value = @'
2024-10-01 20:35:43
'@
$ws.Range('B248').Value = @'
check_availability
'@
$ws.Range('C248').Value = @'
https://www.opentable.com/r/sinatra-bar-and-lounge-nashville
'@
$ws.Range('D248').Value = $errorText
$ws.Range('E248').NumberFormat = '@'
$ws.Range('E248').Value = @'
2024-10-01
'@
$ws.Range('E248').Style = 'Normal'
$ws.Range('F248').Value = @'
20:35:43
'@
